$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.971.42"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.653.02"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3908"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9999"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.113"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.882"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "1.653.68"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06970"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.902"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "23.961.29"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.477"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.013"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.444"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.758"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.481"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.835.06"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.037"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08102"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02970"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.749"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2681"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09164"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.420"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6946"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.456"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.081"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08298"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.222"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
